{"js": "// Update the \"three-digit number divided by one-digit number\" answer\n// table: each cell's old \"dividend\u00f7divisor=quotient, remainder\" text is\n// replaced with a newly generated problem/answer pair. Every old value\n// is unique within the document, so a body-wide search-and-replace per\n// pair is safe and unambiguous.\nconst replacements = [\n  [\"651\u00f76=108, 3\", \"374\u00f76=62, 2\"],\n  [\"453\u00f77=64, 5\", \"332\u00f72=166, 0\"],\n  [\"434\u00f73=144, 2\", \"511\u00f75=102, 1\"],\n  [\"486\u00f77=69, 3\", \"785\u00f75=157, 0\"],\n  [\"989\u00f76=164, 5\", \"445\u00f75=89, 0\"],\n  [\"441\u00f74=110, 1\", \"513\u00f77=73, 2\"],\n  [\"578\u00f77=82, 4\", \"349\u00f79=38, 7\"],\n  [\"818\u00f73=272, 2\", \"951\u00f74=237, 3\"],\n  [\"682\u00f72=341, 0\", \"536\u00f72=268, 0\"],\n  [\"891\u00f75=178, 1\", \"754\u00f74=188, 2\"],\n  [\"358\u00f75=71, 3\", \"831\u00f73=277, 0\"],\n  [\"595\u00f74=148, 3\", \"619\u00f78=77, 3\"],\n  [\"839\u00f78=104, 7\", \"960\u00f79=106, 6\"],\n  [\"124\u00f75=24, 4\", \"906\u00f76=151, 0\"],\n  [\"710\u00f75=142, 0\", \"592\u00f74=148, 0\"],\n  [\"956\u00f75=191, 1\", \"372\u00f73=124, 0\"],\n  [\"858\u00f73=286, 0\", \"321\u00f79=35, 6\"],\n  [\"114\u00f77=16, 2\", \"238\u00f74=59, 2\"],\n  [\"470\u00f72=235, 0\", \"803\u00f77=114, 5\"],\n  [\"574\u00f74=143, 2\", \"606\u00f75=121, 1\"],\n  [\"721\u00f73=240, 1\", \"627\u00f74=156, 3\"],\n  [\"236\u00f73=78, 2\", \"522\u00f78=65, 2\"],\n  [\"292\u00f76=48, 4\", \"905\u00f73=301, 2\"],\n  [\"354\u00f74=88, 2\", \"581\u00f73=193, 2\"],\n  [\"524\u00f75=104, 4\", \"583\u00f78=72, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"three-digit number divided by one-digit number\" answer\n# table: each cell's old \"dividend\u00f7divisor=quotient, remainder\" text is\n# replaced with a newly generated problem/answer pair. Every old value\n# is unique within the document, so Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"651\u00f76=108, 3\"; New = \"374\u00f76=62, 2\" },\n    @{ Old = \"453\u00f77=64, 5\"; New = \"332\u00f72=166, 0\" },\n    @{ Old = \"434\u00f73=144, 2\"; New = \"511\u00f75=102, 1\" },\n    @{ Old = \"486\u00f77=69, 3\"; New = \"785\u00f75=157, 0\" },\n    @{ Old = \"989\u00f76=164, 5\"; New = \"445\u00f75=89, 0\" },\n    @{ Old = \"441\u00f74=110, 1\"; New = \"513\u00f77=73, 2\" },\n    @{ Old = \"578\u00f77=82, 4\"; New = \"349\u00f79=38, 7\" },\n    @{ Old = \"818\u00f73=272, 2\"; New = \"951\u00f74=237, 3\" },\n    @{ Old = \"682\u00f72=341, 0\"; New = \"536\u00f72=268, 0\" },\n    @{ Old = \"891\u00f75=178, 1\"; New = \"754\u00f74=188, 2\" },\n    @{ Old = \"358\u00f75=71, 3\"; New = \"831\u00f73=277, 0\" },\n    @{ Old = \"595\u00f74=148, 3\"; New = \"619\u00f78=77, 3\" },\n    @{ Old = \"839\u00f78=104, 7\"; New = \"960\u00f79=106, 6\" },\n    @{ Old = \"124\u00f75=24, 4\"; New = \"906\u00f76=151, 0\" },\n    @{ Old = \"710\u00f75=142, 0\"; New = \"592\u00f74=148, 0\" },\n    @{ Old = \"956\u00f75=191, 1\"; New = \"372\u00f73=124, 0\" },\n    @{ Old = \"858\u00f73=286, 0\"; New = \"321\u00f79=35, 6\" },\n    @{ Old = \"114\u00f77=16, 2\"; New = \"238\u00f74=59, 2\" },\n    @{ Old = \"470\u00f72=235, 0\"; New = \"803\u00f77=114, 5\" },\n    @{ Old = \"574\u00f74=143, 2\"; New = \"606\u00f75=121, 1\" },\n    @{ Old = \"721\u00f73=240, 1\"; New = \"627\u00f74=156, 3\" },\n    @{ Old = \"236\u00f73=78, 2\"; New = \"522\u00f78=65, 2\" },\n    @{ Old = \"292\u00f76=48, 4\"; New = \"905\u00f73=301, 2\" },\n    @{ Old = \"354\u00f74=88, 2\"; New = \"581\u00f73=193, 2\" },\n    @{ Old = \"524\u00f75=104, 4\"; New = \"583\u00f78=72, 7\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
